$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column AD ("area") rows 2-31 currently hold 8; update them to 10
$ws.Range("AD2:AD31").Value = 10
